$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Move the two "signature" rows (previously 22 and 23) down to the new
#    bottom of the sheet (47 and 48), since 25 new data rows are being
#    inserted into the table above them.
# ---------------------------------------------------------------------------
$ws.Range("B22:C22").Copy($ws.Range("B47:C47"))
$ws.Range("H22:J22").Copy($ws.Range("H47:J47"))
$ws.Range("B23:C23").Copy($ws.Range("B48:C48"))
$ws.Range("H23:J23").Copy($ws.Range("H48:J48"))

# The old locations (rows 22/23) become ordinary data rows, so break their
# merges and clear them out before refilling with data below.
$ws.Range("B22:C22").UnMerge()
$ws.Range("H22:J22").UnMerge()
$ws.Range("B23:C23").UnMerge()
$ws.Range("H23:J23").UnMerge()
$ws.Range("B22:J23").ClearContents()

# ---------------------------------------------------------------------------
# 2. Propagate the "middle of table" row style (currently on row 16) down
#    through row 41, and the "bottom of table" row style (currently on row
#    17) onto the new last data row (42).
# ---------------------------------------------------------------------------
$ws.Range("B17:J17").Copy($ws.Range("B42:J42"))

$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))
$ws.Range("B16:J16").Copy($ws.Range("B19:J19"))
$ws.Range("B16:J16").Copy($ws.Range("B20:J20"))
$ws.Range("B16:J16").Copy($ws.Range("B21:J21"))
$ws.Range("B16:J16").Copy($ws.Range("B22:J22"))
$ws.Range("B16:J16").Copy($ws.Range("B23:J23"))
$ws.Range("B16:J16").Copy($ws.Range("B24:J24"))
$ws.Range("B16:J16").Copy($ws.Range("B25:J25"))
$ws.Range("B16:J16").Copy($ws.Range("B26:J26"))
$ws.Range("B16:J16").Copy($ws.Range("B27:J27"))
$ws.Range("B16:J16").Copy($ws.Range("B28:J28"))
$ws.Range("B16:J16").Copy($ws.Range("B29:J29"))
$ws.Range("B16:J16").Copy($ws.Range("B30:J30"))
$ws.Range("B16:J16").Copy($ws.Range("B31:J31"))
$ws.Range("B16:J16").Copy($ws.Range("B32:J32"))
$ws.Range("B16:J16").Copy($ws.Range("B33:J33"))
$ws.Range("B16:J16").Copy($ws.Range("B34:J34"))
$ws.Range("B16:J16").Copy($ws.Range("B35:J35"))
$ws.Range("B16:J16").Copy($ws.Range("B36:J36"))
$ws.Range("B16:J16").Copy($ws.Range("B37:J37"))
$ws.Range("B16:J16").Copy($ws.Range("B38:J38"))
$ws.Range("B16:J16").Copy($ws.Range("B39:J39"))
$ws.Range("B16:J16").Copy($ws.Range("B40:J40"))
$ws.Range("B16:J16").Copy($ws.Range("B41:J41"))

# ---------------------------------------------------------------------------
# 3. Fill in the employee / period overdue-balance table (rows 16-42).
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1044923865"
$ws.Range("D16").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E16").Value = "2111"
$ws.Range("F16").Value = 35112
$ws.Range("G16").Value = 877803
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1044923865"
$ws.Range("D17").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E17").Value = "2110"
$ws.Range("F17").Value = 35112
$ws.Range("G17").Value = 877803
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1044923865"
$ws.Range("D18").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E18").Value = "2109"
$ws.Range("F18").Value = 35112
$ws.Range("G18").Value = 877803
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1044923865"
$ws.Range("D19").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E19").Value = "2108"
$ws.Range("F19").Value = 35112
$ws.Range("G19").Value = 877803
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1044923865"
$ws.Range("D20").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E20").Value = "2107"
$ws.Range("F20").Value = 35112
$ws.Range("G20").Value = 877803
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1044923865"
$ws.Range("D21").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E21").Value = "2106"
$ws.Range("F21").Value = 35112
$ws.Range("G21").Value = 877803
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1044923865"
$ws.Range("D22").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E22").Value = "2105"
$ws.Range("F22").Value = 35112
$ws.Range("G22").Value = 877803
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1044923865"
$ws.Range("D23").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E23").Value = "2104"
$ws.Range("F23").Value = 35112
$ws.Range("G23").Value = 877803
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1044923865"
$ws.Range("D24").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E24").Value = "2103"
$ws.Range("F24").Value = 35112
$ws.Range("G24").Value = 877803
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1044923865"
$ws.Range("D25").Value = "NIBER ROBERTO FABREGA MENDOZA"
$ws.Range("E25").Value = "2102"
$ws.Range("F25").Value = 35112
$ws.Range("G25").Value = 877803
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "3830224"
$ws.Range("D26").Value = "JESUS ALFARO OSPINO"
$ws.Range("E26").Value = "2111"
$ws.Range("F26").Value = 35112
$ws.Range("G26").Value = 908526
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "3830224"
$ws.Range("D27").Value = "JESUS ALFARO OSPINO"
$ws.Range("E27").Value = "2110"
$ws.Range("F27").Value = 35112
$ws.Range("G27").Value = 908526
$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "3830224"
$ws.Range("D28").Value = "JESUS ALFARO OSPINO"
$ws.Range("E28").Value = "2109"
$ws.Range("F28").Value = 35112
$ws.Range("G28").Value = 908526
$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "3830224"
$ws.Range("D29").Value = "JESUS ALFARO OSPINO"
$ws.Range("E29").Value = "2108"
$ws.Range("F29").Value = 35112
$ws.Range("G29").Value = 908526
$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "3830224"
$ws.Range("D30").Value = "JESUS ALFARO OSPINO"
$ws.Range("E30").Value = "2107"
$ws.Range("F30").Value = 35112
$ws.Range("G30").Value = 908526
$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "3830224"
$ws.Range("D31").Value = "JESUS ALFARO OSPINO"
$ws.Range("E31").Value = "2106"
$ws.Range("F31").Value = 35112
$ws.Range("G31").Value = 908526
$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "3830224"
$ws.Range("D32").Value = "JESUS ALFARO OSPINO"
$ws.Range("E32").Value = "2105"
$ws.Range("F32").Value = 35112
$ws.Range("G32").Value = 908526
$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "3830224"
$ws.Range("D33").Value = "JESUS ALFARO OSPINO"
$ws.Range("E33").Value = "2104"
$ws.Range("F33").Value = 35112
$ws.Range("G33").Value = 908526
$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "3830224"
$ws.Range("D34").Value = "JESUS ALFARO OSPINO"
$ws.Range("E34").Value = "2103"
$ws.Range("F34").Value = 35112
$ws.Range("G34").Value = 908526
$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "3830224"
$ws.Range("D35").Value = "JESUS ALFARO OSPINO"
$ws.Range("E35").Value = "2102"
$ws.Range("F35").Value = 35112
$ws.Range("G35").Value = 908526
$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "1007954996"
$ws.Range("D36").Value = "ADRIAN ANDRES OSPINO PEREZ"
$ws.Range("E36").Value = "2304"
$ws.Range("F36").Value = 46400
$ws.Range("G36").Value = 1160000
$ws.Range("B37").Value = "CC"
$ws.Range("C37").Value = "1044920005"
$ws.Range("D37").Value = "JAIDER JOSE ROMERO JULIO"
$ws.Range("E37").Value = "2507"
$ws.Range("F37").Value = 52000
$ws.Range("G37").Value = 781242
$ws.Range("B38").Value = "CC"
$ws.Range("C38").Value = "1044920005"
$ws.Range("D38").Value = "JAIDER JOSE ROMERO JULIO"
$ws.Range("E38").Value = "2506"
$ws.Range("F38").Value = 52000
$ws.Range("G38").Value = 781242
$ws.Range("B39").Value = "CC"
$ws.Range("C39").Value = "1001834677"
$ws.Range("D39").Value = "XAVIER HUMBERTO CASTELLON MELENDREZ"
$ws.Range("E39").Value = "2009"
$ws.Range("F39").Value = 35112
$ws.Range("G39").Value = 877803
$ws.Range("B40").Value = "CC"
$ws.Range("C40").Value = "143357561"
$ws.Range("D40").Value = "HARRY JOSE CARABALLO CASTRO"
$ws.Range("E40").Value = "2012"
$ws.Range("F40").Value = 35112
$ws.Range("G40").Value = 877803
$ws.Range("B41").Value = "CC"
$ws.Range("C41").Value = "143357561"
$ws.Range("D41").Value = "HARRY JOSE CARABALLO CASTRO"
$ws.Range("E41").Value = "2011"
$ws.Range("F41").Value = 19897
$ws.Range("G41").Value = 877803
$ws.Range("B42").Value = "CC"
$ws.Range("C42").Value = "1193368427"
$ws.Range("D42").Value = "YOINER VILLAMIZAR MEJIA"
$ws.Range("E42").Value = "2306"
$ws.Range("F42").Value = 46400
$ws.Range("G42").Value = 1160000

# ---------------------------------------------------------------------------
# 4. Misc header/summary cell value updates.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 989161
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 17

# ---------------------------------------------------------------------------
# 5. Column D widened to fit the longer employee names now in the table.
# ---------------------------------------------------------------------------
$ws.Range("D1").EntireColumn.ColumnWidth = 41
